# Generate Report for Handback
# Updates the "generate date" timestamps recorded on the handback-status
# report: the Overview sheet's "Latest HO Xliff Generate Date" column, and
# the per-locale sheets' "Correspond Handoff/Handback Datetime" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-04 01:11:33"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-04 01:11:28"
$zhcn.Range("K2").Value = "2016-09-04 01:11:55"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-04 01:11:33"
$dede.Range("K2").Value = "2016-09-04 01:12:04"
